$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.459.32'
$ws.Range('E2').Value = '  +2.67%  '
$ws.Range('D3').Value = '3.782.07'
$ws.Range('E3').Value = '  +0.25%  '
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.64%  '
$ws.Range('D5').Value = '626.34'
$ws.Range('E5').Value = '  +4.68%  '
$ws.Range('D6').Value = '166.70'
$ws.Range('E6').Value = '  +2.86%  '
$ws.Range('D7').Value = '3.780.89'
$ws.Range('E7').Value = '  +0.37%  '
$ws.Range('E8').Value = '  -0.24%  '
$ws.Range('E9').Value = '  +2.14%  '
$ws.Range('E10').Value = '  +3.23%  '
$ws.Range('D11').Value = '0.460'
$ws.Range('E11').Value = '  +3.83%  '
$ws.Range('D12').Value = '6.75'
$ws.Range('E12').Value = '  +2.54%  '
$ws.Range('E13').Value = '  +2.28%  '
$ws.Range('D14').Value = '35.99'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('D15').Value = '4.413.76'
$ws.Range('E15').Value = '  +0.13%  '
$ws.Range('D16').Value = '3.781.02'
$ws.Range('E16').Value = '  -0.30%  '
$ws.Range('D17').Value = '69.400.01'
$ws.Range('E17').Value = '  +2.55%  '
$ws.Range('D18').Value = '17.74'
$ws.Range('E18').Value = '  -1.96%  '
$ws.Range('D19').Value = '7.11'
$ws.Range('E19').Value = '  +1.95%  '
$ws.Range('E20').Value = '  -0.86%  '
$ws.Range('D21').Value = '470.29'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').Value = '9.61'
$ws.Range('E22').Value = '  +2.24%  '
$ws.Range('E23').Value = '  +3.14%  '
$ws.Range('D24').Value = '0.0000149'
$ws.Range('E24').Value = '  +6.09%  '
$ws.Range('D25').Value = '83.39'
$ws.Range('E25').Value = '  +0.73%  '
$ws.Range('D26').Value = '12.11'
$ws.Range('E26').Value = '  +2.71%  '
$ws.Range('E27').Value = '  +4.99%  '
$ws.Range('D28').Value = '10.05'
$ws.Range('E28').Value = '  +2.07%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D30').Value = '3.930.00'
$ws.Range('E30').Value = '  +0.05%  '
$ws.Range('E31').Value = '  +4.19%  '
$ws.Range('D32').Value = '2.26'
$ws.Range('E32').Value = '  +3.71%  '
$ws.Range('D33').Value = '7.23'
$ws.Range('E33').Value = '  +0.52%  '
$ws.Range('D34').Value = '28.91'
$ws.Range('E34').Value = '  +0.65%  '
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.731.64'
$ws.Range('E36').Value = '  +0.21%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').Value = '9.00'
$ws.Range('E37').Value = '  +1.23%  '
$ws.Range('E38').Value = '  +13.88%  '
$ws.Range('E39').Value = '  +3.10%  '
$ws.Range('D40').Value = '3.43'
$ws.Range('E40').Value = '  +9.53%  '
$ws.Range('D41').Value = '5.84'
$ws.Range('E41').Value = '  +1.34%  '
$ws.Range('D42').Value = '0.968'
$ws.Range('E42').Value = '  -0.70%  '
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('E45').Value = '  +2.07%  '
$ws.Range('B46').Value = 'Arweave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D46').Value = '43.21'
$ws.Range('E46').Value = '  -0.15%  '
$ws.Range('B47').Value = 'Monero'
$ws.Range('C47').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D47').Value = '153.00'
$ws.Range('E47').Value = '  +1.51%  '
$ws.Range('B48').Value = 'OKB'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D48').Value = '46.78'
$ws.Range('E48').Value = '  -0.70%  '
$ws.Range('B49').Value = 'Stacks'
$ws.Range('C49').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D49').Value = '1.92'
$ws.Range('E49').Value = '  +5.15%  '
$ws.Range('D50').Value = '8.43'
$ws.Range('E50').Value = '  +2.31%  '
$ws.Range('D51').Value = '1.36'
$ws.Range('E51').Value = '  +0.98%  '
